$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.708.98'
$ws.Range("E2").Value = '  +3.69%  '
$ws.Range("D3").Value = '2.346.28'
$ws.Range("E3").Value = '  +2.95%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '517.35'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.69%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '133.61'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.01%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.43%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.535'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.18%  '
$ws.Range("D9").Value = '2.343.03'
$ws.Range("E9").Value = '  +2.11%  '
$ws.Range("E10").Value = '  +6.53%  '
$ws.Range("E11").Value = '  +0.33%  '
$ws.Range("E12").Value = '  +7.19%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.340'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.88%  '
$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.74'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.48%  '
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = '2.764.25'
$ws.Range("E15").Value = '  +2.91%  '
$ws.Range("D16").Value = '56.713.41'
$ws.Range("E16").Value = '  +3.56%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000134'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.15%  '
$ws.Range("D18").Value = '2.330.52'
$ws.Range("E18").Value = '  +1.86%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.40'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.72%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.24'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.63%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '320.37'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.41%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.58'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.53%  '
$ws.Range("E23").Value = '  +0.27%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '60.76'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.01%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.999'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.00%  '
$ws.Range("E26").Value = '  +5.90%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.72'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.03%  '
$ws.Range("B28").Value = 'Monero'
$ws.Range("C28").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '170.97'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.08%  '
$ws.Range("B29").Value = 'Fetch.AI'
$ws.Range("C29").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.22'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +9.48%  '
$ws.Range("D30").Value = '0.0₃0736'
$ws.Range("E30").Value = '  +5.06%  '
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.67'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.71%  '
$ws.Range("B32").Value = 'Aptos'
$ws.Range("C32").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.21'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.58%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.25'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.84%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.49%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.940'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.03%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.24'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.60%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.97'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.86%  '
$ws.Range("E39").Value = '  +7.91%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '37.50'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.89%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.379'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.98%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '138.01'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +9.21%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.56'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.98%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '275.72'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +9.27%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.04'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.34%  '
$ws.Range("E46").Value = '  +3.23%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0502'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.04%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.558'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.71%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0216'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.71%  '
$ws.Range("E50").Value = '  +1.10%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '16.71'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.32%  '
